# Update report xlsx file
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend the title merge from A1:C2 to A1:D2 (re-merging the superset range
# absorbs / replaces the existing A1:C2 merge in one step)
$ws.Range("A1:D2").MergeCells = $true

# Fill in the new "remark" row (row 19) with the reported tasks
$ws.Range("B19").Value = "메일 서버 구축(postfix)`n메일 서버 작동 원리 공부"
$ws.Range("B19").WrapText = $true
$ws.Range("C19").Value = "spamassassin 적용"

# Update the current selection to match the new merged title range
$ws.Range("A1:D2").Select()

# Resize the workbook window (as captured in the saved file)
$excel.Width = 23040
$excel.Height = 8976
